$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend formatting (styles) down from row 178 into the new rows ---
# Rows 179:192 get the full A:D format (text date, counts, percent formula cell);
# row 193 (the new last row) only gets A:C, matching the source (no D formula there).
$ws.Range("A178:D178").Copy()
$ws.Range("A179:D192").PasteSpecial(-4122)
$ws.Range("A178:C178").Copy()
$ws.Range("A193:C193").PasteSpecial(-4122)

# --- Fill in the new data rows (dates 2020-10-01 .. 2020-10-15) ---
$ws.Range("A179").Value = "2020-10-01"
$ws.Range("B179").Value = 65
$ws.Range("C179").Value = 62

$ws.Range("A180").Value = "2020-10-02"
$ws.Range("B180").Value = 67
$ws.Range("C180").Value = 64

$ws.Range("A181").Value = "2020-10-03"
$ws.Range("B181").Value = 38
$ws.Range("C181").Value = 38

$ws.Range("A182").Value = "2020-10-04"
$ws.Range("B182").Value = 54
$ws.Range("C182").Value = 49

$ws.Range("A183").Value = "2020-10-05"
$ws.Range("B183").Value = 53
$ws.Range("C183").Value = 50

$ws.Range("A184").Value = "2020-10-06"
$ws.Range("B184").Value = 42
$ws.Range("C184").Value = 42

$ws.Range("A185").Value = "2020-10-07"
$ws.Range("B185").Value = 45
$ws.Range("C185").Value = 43

$ws.Range("A186").Value = "2020-10-08"
$ws.Range("B186").Value = 56
$ws.Range("C186").Value = 54

$ws.Range("A187").Value = "2020-10-09"
$ws.Range("B187").Value = 62
$ws.Range("C187").Value = 59

$ws.Range("A188").Value = "2020-10-10"
$ws.Range("B188").Value = 42
$ws.Range("C188").Value = 41

$ws.Range("A189").Value = "2020-10-11"
$ws.Range("B189").Value = 48
$ws.Range("C189").Value = 48

$ws.Range("A190").Value = "2020-10-12"
$ws.Range("B190").Value = 59
$ws.Range("C190").Value = 56

$ws.Range("A191").Value = "2020-10-13"
$ws.Range("B191").Value = 57
$ws.Range("C191").Value = 55

$ws.Range("A192").Value = "2020-10-14"
$ws.Range("B192").Value = 46
$ws.Range("C192").Value = 44

$ws.Range("A193").Value = "2020-10-15"
$ws.Range("B193").Value = 65
$ws.Range("C193").Value = 63

# --- Add the % formula for the new rows in one shot so it forms a shared formula group ---
$ws.Range("D179:D192").Formula = "=C179/B179"

# --- Update the visible selection to match the new bottom of the sheet ---
$ws.Range("F193").Select()